$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 76, shifting existing rows 76:121 down to 77:122
$ws.Rows.Item(76).Insert()

# Populate the new row 76 with the new record
$ws.Range("A76").Value = 3
$ws.Range("B76").Value = "Femacal de La Calera"
$ws.Range("C76").Value = "Coquimbo"
$ws.Range("D76").Value = 44455
$ws.Range("E76").Value = 5
$ws.Range("F76").Value = 100112010
$ws.Range("G76").Value = "Achicoria"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 130
$ws.Range("K76").Value = 6000
$ws.Range("L76").Value = 6500
$ws.Range("M76").Value = 6269
$ws.Range("N76").Value = '$/caja 16 unidades'
$ws.Range("O76").Value = "Provincia de Quillota"
$ws.Range("P76").Value = 392
$ws.Range("Q76").Value = 16
$ws.Range("R76").Value = "Hortaliza"
